$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the counts in row 2 (B2:E2), resetting B2/D2/E2/C2 to the default
# "Normal" style (they currently carry style s="3" on C2/D2/E2; B2 already
# had no explicit style) so that no "s" attribute remains on any of them.
$ws.Range("B2:E2").Style = "Normal"

$ws.Range("B2").Value = 1177
$ws.Range("C2").Value = 1180
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 5

# New data point added at G16 (default/no explicit style, overriding the
# column's inherited style).
$ws.Range("G16").Style = "Normal"
$ws.Range("G16").Value = 1182

# Update the selected cell shown in the sheet view to G16.
$ws.Range("G16").Select()
